# Eurostat domestic-percentage workbook refresh
# Applies the "updated data sets to recent" commit:
#  - refreshes a handful of 2021-07/08/09 (T/U/V) data points with newer
#    figures (and re-stripes their banding style where the refreshed row
#    shifted parity)
#  - extends the United Kingdom row (row 36) so columns N:V now show the
#    ":" (not-available) placeholder like the rest of the sheet
#  - re-zooms the sheet and freezes column A while moving the selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: stamp a target range with the same direct formatting (number
# format + fill, i.e. the workbook's "style" index) as a known reference
# cell, without touching any values already in the target range.
# ---------------------------------------------------------------------
function Copy-CellStyle {
    param(
        [string]$RefCell,
        [string]$TargetRange
    )
    $ws.Range($RefCell).Copy() | Out-Null
    $ws.Range($TargetRange).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = $false
}

# Stable reference cells (outside every range touched below) for each of
# the four direct-format "styles" used across the data grid:
#   B2 -> shaded numeric   (style 3)
#   K9 -> shaded ":" text  (style 4)
#   B3 -> plain numeric    (style 5)
#   C3 -> plain ":" text   (style 6)
$styleShadedNum  = "B2"
$styleShadedText = "K9"
$stylePlainNum   = "B3"
$stylePlainText  = "C3"

# ---------------------------------------------------------------------
# Row 2 - European Union (27 countries): revised Aug value, new Sep value
# ---------------------------------------------------------------------
$ws.Range("U2").Value = -4.1900000000000004
Copy-CellStyle $styleShadedNum "V2"
$ws.Range("V2").Value = 15.38

# Row 9 - Germany: new Sep value
Copy-CellStyle $styleShadedNum "V9"
$ws.Range("V9").Value = 6.63

# Row 10 - Belgium: new Sep value
Copy-CellStyle $stylePlainNum "V10"
$ws.Range("V10").Value = 29.71

# Row 12 - Czechia: new Sep value
Copy-CellStyle $stylePlainNum "V12"
$ws.Range("V12").Value = -8.94

# Row 16 - Estonia: revised Aug value, new Sep value
$ws.Range("U16").Value = -0.38
Copy-CellStyle $stylePlainNum "V16"
$ws.Range("V16").Value = 10.74

# Row 18 - Spain: restripe only (values unchanged)
Copy-CellStyle $styleShadedNum "T18:V18"

# Row 19 - France: restripe + revised Jul/Aug/Sep values
Copy-CellStyle $stylePlainNum "T19:V19"
$ws.Range("T19").Value = 24.59
$ws.Range("U19").Value = 13.28
$ws.Range("V19").Value = 1.92

# Row 20 - Croatia: restripe Jul/Aug + revised Sep value
Copy-CellStyle $styleShadedNum "T20:U20"
Copy-CellStyle $styleShadedNum "V20"
$ws.Range("V20").Value = -32.58

# Row 21 - Italy: restripe only (values unchanged)
Copy-CellStyle $stylePlainNum "T21:V21"

# Row 22 - Cyprus: restripe + revised Aug/Sep values
Copy-CellStyle $styleShadedNum "T22"
Copy-CellStyle $styleShadedNum "U22"
$ws.Range("U22").Value = 97.16
Copy-CellStyle $styleShadedNum "V22"
$ws.Range("V22").Value = 110.22

# Row 23 - Latvia: restripe Jul/Aug + new Sep value
Copy-CellStyle $stylePlainNum "T23:U23"
Copy-CellStyle $stylePlainNum "V23"
$ws.Range("V23").Value = 29.57

# Row 24 - Lithuania: restripe only (values unchanged)
Copy-CellStyle $styleShadedNum "T24:V24"

# Row 25 - Luxembourg: restripe only (values unchanged)
Copy-CellStyle $stylePlainNum "T25:V25"

# Row 26 - Hungary: restripe only (values unchanged)
Copy-CellStyle $styleShadedNum "T26:V26"

# Row 27 - Malta: restripe Jul/Aug + new Sep value
Copy-CellStyle $stylePlainNum "T27:U27"
Copy-CellStyle $stylePlainNum "V27"
$ws.Range("V27").Value = -15.9

# Row 28 - Netherlands: restripe Jul/Aug + new Sep value
Copy-CellStyle $styleShadedNum "T28:U28"
Copy-CellStyle $styleShadedNum "V28"
$ws.Range("V28").Value = 154.94999999999999

# Row 29 - Austria: restripe Jul/Aug + new Sep value
Copy-CellStyle $stylePlainNum "T29:U29"
Copy-CellStyle $stylePlainNum "V29"
$ws.Range("V29").Value = -5.33

# Row 30 - Poland: restripe only (values unchanged)
Copy-CellStyle $styleShadedNum "T30:V30"

# Row 31 - Portugal: restripe only (values unchanged)
Copy-CellStyle $stylePlainNum "T31:V31"

# Row 32 - Romania: restripe Jul, re-shade still-missing Aug/Sep
Copy-CellStyle $styleShadedNum "T32"
Copy-CellStyle $styleShadedText "U32:V32"
$ws.Range("U32").Value = ":"
$ws.Range("V32").Value = ":"

# Row 33 - Slovenia: restripe Jul/Aug + new Sep value
Copy-CellStyle $stylePlainNum "T33:U33"
Copy-CellStyle $stylePlainNum "V33"
$ws.Range("V33").Value = 50.31

# Row 34 - Slovakia: restripe Jul/Aug + new Sep value
Copy-CellStyle $styleShadedNum "T34:U34"
Copy-CellStyle $styleShadedNum "V34"
$ws.Range("V34").Value = 7.91

# ---------------------------------------------------------------------
# Row 36 - United Kingdom: no data reported from 2021-01 (N) onward, so
# fill N36:V36 with the same ":" placeholder/style used for I36:M36
# ---------------------------------------------------------------------
Copy-CellStyle "M36" "N36:V36"
$ws.Range("N36:V36").Value = ":"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet view: zoom to 90%, freeze column A, move the active selection
# ---------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 90
$ws.Range("B1").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("X13").Select() | Out-Null
